$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from C1 (header style) to new header cells D1 and E1
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update header row (row 1) text
$ws.Range("C1").Value = "Frecuencia del primer armonico"
$ws.Range("D1").Value = "Frecuencia del segundo armonico"
$ws.Range("E1").Value = "Frecuencia tercer armonico"

# Update data rows 2-58 with new C values and new D/E columns (frequencies of 1st/2nd/3rd harmonics)
$ws.Range("C2").Value = 318.3664507629819
$ws.Range("D2").Value = 159.4027884509824
$ws.Range("E2").Value = 476.451860797014
$ws.Range("C3").Value = 375.5353777625487
$ws.Range("D3").Value = 187.7676888812739
$ws.Range("E3").Value = 564.67363371595
$ws.Range("C4").Value = 425.2604720391246
$ws.Range("D4").Value = 284.0739953221346
$ws.Range("E4").Value = 142.0369976610673
$ws.Range("C5").Value = 332.7472934209309
$ws.Range("D5").Value = 166.1885814749703
$ws.Range("E5").Value = 497.0852225409462
$ws.Range("C6").Value = 328.8391066860258
$ws.Range("D6").Value = 492.0141049574781
$ws.Range("E6").Value = 164.0047016524923
$ws.Range("C7").Value = 376.3616090573805
$ws.Range("D7").Value = 187.6661806329876
$ws.Range("E7").Value = 554.7645595677159
$ws.Range("C8").Value = 369.4012033234649
$ws.Range("D8").Value = 185.2736128354509
$ws.Range("E8").Value = 556.5848534046418
$ws.Range("C9").Value = 301.875
$ws.Range("D9").Value = 451.875
$ws.Range("E9").Value = 151.25
$ws.Range("C10").Value = 411.3360323886636
$ws.Range("D10").Value = 205.3441295546563
$ws.Range("E10").Value = 2267.206477732793
$ws.Range("C11").Value = 285.1772960781864
$ws.Range("D11").Value = 569.8534018293444
$ws.Range("E11").Value = 133.8178173161255
$ws.Range("C12").Value = 390.948483389504
$ws.Range("D12").Value = 195.7952174610809
$ws.Range("E12").Value = 577.1144278606962
$ws.Range("C13").Value = 504.6771617613504
$ws.Range("D13").Value = 336.7556468172484
$ws.Range("E13").Value = 172.4845995893229
$ws.Range("C14").Value = 187.2981700753503
$ws.Range("D14").Value = 374.5963401506997
$ws.Range("E14").Value = 561.176892716182
$ws.Range("C15").Value = 317.4217087260567
$ws.Range("D15").Value = 159.3189419276378
$ws.Range("E15").Value = 477.1460423634339
$ws.Range("C16").Value = 367.3161804284637
$ws.Range("D16").Value = 184.7346323608572
$ws.Range("E16").Value = 550.758962213371
$ws.Range("C17").Value = 370.5844572896594
$ws.Range("D17").Value = 184.9710982658962
$ws.Range("E17").Value = 555.5555555555557
$ws.Range("C18").Value = 365.5549274543055
$ws.Range("D18").Value = 183.9080459770121
$ws.Range("E18").Value = 540.4183154324483
$ws.Range("C19").Value = 321.2598425196848
$ws.Range("D19").Value = 481.8897637795276
$ws.Range("E19").Value = 160.6299212598424
$ws.Range("C20").Value = 352.8624037495815
$ws.Range("D20").Value = 176.7659859390692
$ws.Range("E20").Value = 530.297957817208
$ws.Range("C21").Value = 311.4511352418558
$ws.Range("D21").Value = 465.942744323791
$ws.Range("E21").Value = 163.376110562685
$ws.Range("C22").Value = 498.5190248348144
$ws.Range("D22").Value = 169.5146958304849
$ws.Range("E22").Value = 332.6498063340168
$ws.Range("C23").Value = 434.164375770215
$ws.Range("D23").Value = 652.5737036685941
$ws.Range("E23").Value = 216.513413593706
$ws.Range("C24").Value = 510.1243339253997
$ws.Range("D24").Value = 247.2468916518646
$ws.Range("E24").Value = 625.9325044404977
$ws.Range("C25").Value = 450.8231321232588
$ws.Range("D25").Value = 299.8251221130072
$ws.Range("E25").Value = 149.7919556171983
$ws.Range("C26").Value = 1521.54026167429
$ws.Range("D26").Value = 509.3075204765451
$ws.Range("E26").Value = 331.0286139772365
$ws.Range("C27").Value = 465.1162790697672
$ws.Range("D27").Value = 311.3255093959242
$ws.Range("E27").Value = 154.6547627618975
$ws.Range("C28").Value = 583.8128116609123
$ws.Range("D28").Value = 402.7617951668581
$ws.Range("E28").Value = 194.0928270042195
$ws.Range("C29").Value = 610.3067733117346
$ws.Range("D29").Value = 458.843778475245
$ws.Range("E29").Value = 1678.647362559482
$ws.Range("C30").Value = 467.2701949860721
$ws.Range("D30").Value = 623.2590529247909
$ws.Range("E30").Value = 1866.991643454039
$ws.Range("C31").Value = 521.3818860877682
$ws.Range("D31").Value = 652.1008403361348
$ws.Range("E31").Value = 391.4098972922502
$ws.Range("C32").Value = 454.5985158008843
$ws.Range("D32").Value = 605.3531226548821
$ws.Range("E32").Value = 303.1768531643456
$ws.Range("C33").Value = 481.6985513656973
$ws.Range("D33").Value = 241.0704412252571
$ws.Range("E33").Value = 720.5573371668697
$ws.Range("C34").Value = 533.4875650665126
$ws.Range("D34").Value = 266.0497397339504
$ws.Range("E34").Value = 3190.283400809717
$ws.Range("C35").Value = 530.3494788473331
$ws.Range("D35").Value = 264.8681790312694
$ws.Range("E35").Value = 3181.483752299203
$ws.Range("C36").Value = 582.8592402275644
$ws.Range("D36").Value = 292.1636997614241
$ws.Range("E36").Value = 1456.41402092127
$ws.Range("C37").Value = 549.7342444950646
$ws.Range("D37").Value = 275.5179520555375
$ws.Range("E37").Value = 822.2149907799112
$ws.Range("C38").Value = 492.5081433224755
$ws.Range("D38").Value = 245.8197611292071
$ws.Range("E38").Value = 738.3279044516826
$ws.Range("C39").Value = 501.5027908973807
$ws.Range("D39").Value = 250.7513954486903
$ws.Range("E39").Value = 752.2541863460715
$ws.Range("C40").Value = 516.3488843813384
$ws.Range("D40").Value = 258.0121703853956
$ws.Range("E40").Value = 3088.68154158215
$ws.Range("C41").Value = 380.6271902237395
$ws.Range("D41").Value = 183.6643004762332
$ws.Range("E41").Value = 571.1204960014375
$ws.Range("C42").Value = 461.8342191563952
$ws.Range("D42").Value = 231.1432771683817
$ws.Range("E42").Value = 2524.482641637454
$ws.Range("C43").Value = 458.9091503802438
$ws.Range("D43").Value = 229.9452122005077
$ws.Range("E43").Value = 2516.313680595306
$ws.Range("C44").Value = 483.5367257655998
$ws.Range("D44").Value = 243.1498963849872
$ws.Range("E44").Value = 724.8445774810034
$ws.Range("C45").Value = 486.2393236131484
$ws.Range("D45").Value = 243.9824001380384
$ws.Range("E45").Value = 727.8060564230868
$ws.Range("C46").Value = 491.3085004775548
$ws.Range("D46").Value = 246.0362941738304
$ws.Range("E46").Value = 739.6370582617001
$ws.Range("C47").Value = 194.1899005654122
$ws.Range("D47").Value = 373.5620978748293
$ws.Range("E47").Value = 573.2111522713985
$ws.Range("C48").Value = 274.3333699111163
$ws.Range("D48").Value = 548.2278064303737
$ws.Range("E48").Value = 823.0001097333479
$ws.Range("C49").Value = 544.9951409135083
$ws.Range("D49").Value = 272.8862973760929
$ws.Range("E49").Value = 817.8814382896012
$ws.Range("C50").Value = 601.8423746161716
$ws.Range("D50").Value = 301.107285754164
$ws.Range("E50").Value = 906.2994323997391
$ws.Range("C51").Value = 551.0264802142219
$ws.Range("D51").Value = 276.7033620946145
$ws.Range("E51").Value = 827.7298423088369
$ws.Range("C52").Value = 524.6837249120135
$ws.Range("D52").Value = 174.6409207647675
$ws.Range("E52").Value = 350.4232854561023
$ws.Range("C53").Value = 438.8312105319492
$ws.Range("D53").Value = 293.2676870384248
$ws.Range("E53").Value = 146.8479075243495
$ws.Range("C54").Value = 392.8244074898521
$ws.Range("D54").Value = 197.4597354982316
$ws.Range("E54").Value = 587.1415477281653
$ws.Range("C55").Value = 388.1645717854635
$ws.Range("D55").Value = 194.5507065344682
$ws.Range("E55").Value = 582.4029978921071
$ws.Range("C56").Value = 203.3640596787991
$ws.Range("D56").Value = 406.3479996198803
$ws.Range("E56").Value = 608.5717000855266
$ws.Range("C57").Value = 458.5987261146497
$ws.Range("D57").Value = 229.2993630573251
$ws.Range("E57").Value = 687.8980891719748
$ws.Range("C58").Value = 499.5222344726199
$ws.Range("D58").Value = 251.3781697905179
$ws.Range("E58").Value = 754.1345093715545
